# Add translation settings to the ODK-X "settings" sheet:
#  - new title-translation columns (D/E) and a locale-translation block (F/G/H)
#  - a "display.locale.text" setting (row 1, col F) describing the available locales
#  - three new rows (7-9) enumerating the "default" / "pt" / "sw" locales
#  - make "settings" the active sheet/tab

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# --- Row 1 (headers) -------------------------------------------------
$ws.Range("D1").Value = "display.title.text.pt"
$ws.Range("E1").Value = "display.title.text.sw"
$ws.Range("F1").Value = "display.locale.text"
$ws.Range("G1").Value = "display.locale.text.pt"
$ws.Range("H1").Value = "display.locale.text.sw"

# --- Row 5 (form title translations) ----------------------------------
$ws.Range("D5").Value = "Household Members Absent"
$ws.Range("E5").Value = "Household Members Absent"

# --- Rows 7-9 (locale definitions) ------------------------------------
$ws.Range("A7").Value = "default"
$ws.Range("F7").Value = "English"
$ws.Range("G7").Value = "English"
$ws.Range("H7").Value = "English"

$ws.Range("A8").Value = "pt"
$ws.Range("F8").Value = "Português"
$ws.Range("G8").Value = "Português"
$ws.Range("H8").Value = "Português"

$ws.Range("A9").Value = "sw"
$ws.Range("F9").Value = "Kiswahili"
$ws.Range("G9").Value = "Kiswahili"
$ws.Range("H9").Value = "Kiswahili"

# --- Make "settings" the active sheet/tab -----------------------------
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
